$p = $ppt.ActivePresentation
Write-Output ("SlideMasterCount_before=" + $p.Slides.Count)
$p.Slides.Add(2, 1) | Out-Null
Write-Output ("SlideCount_after=" + $p.Slides.Count)
